$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": update row 7 (Order ID 13, Ketki) ---
$ordersSheet = $wb.Worksheets.Item("All Orders")
$ordersSheet.Range("H7").Value = "CANCELLED"
$ordersSheet.Range("M7").Value = "test order"

# --- Sheet "Daily Summary": update row 4 (2026-01-13 totals) ---
$summarySheet = $wb.Worksheets.Item("Daily Summary")
$summarySheet.Range("D4").Value = 11
$summarySheet.Range("E4").Value = 35
$summarySheet.Range("G4").Value = 35
